# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 2;   I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 28;  I = "ba"; J = "Appreciation" },
    @{ Row = 33;  I = "%";  J = "Uninterpretable" },
    @{ Row = 42;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 48;  I = "%";  J = "Uninterpretable" },
    @{ Row = 64;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 67;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 71;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 91;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 102; I = "sv"; J = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

$wb.Save()
